# Update the five-row division practice table: replace each problem's
# text with the newly generated problem, cell by cell, preserving all
# existing run/paragraph formatting (font, size, justification).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "580÷3="
$t.Cell(1,2).Range.Text = "384÷2="
$t.Cell(1,3).Range.Text = "960÷9="
$t.Cell(1,4).Range.Text = "365÷8="
$t.Cell(1,5).Range.Text = "767÷9="

$t.Cell(5,1).Range.Text = "926÷6="
$t.Cell(5,2).Range.Text = "770÷4="
$t.Cell(5,3).Range.Text = "467÷6="
$t.Cell(5,4).Range.Text = "962÷2="
$t.Cell(5,5).Range.Text = "395÷6="

$t.Cell(9,1).Range.Text = "574÷8="
$t.Cell(9,2).Range.Text = "594÷6="
$t.Cell(9,3).Range.Text = "297÷9="
$t.Cell(9,4).Range.Text = "897÷4="
$t.Cell(9,5).Range.Text = "164÷5="

$t.Cell(13,1).Range.Text = "466÷2="
$t.Cell(13,2).Range.Text = "671÷8="
$t.Cell(13,3).Range.Text = "225÷8="
$t.Cell(13,4).Range.Text = "360÷8="
$t.Cell(13,5).Range.Text = "387÷4="

$t.Cell(17,1).Range.Text = "356÷4="
$t.Cell(17,2).Range.Text = "573÷6="
$t.Cell(17,3).Range.Text = "278÷5="
$t.Cell(17,4).Range.Text = "192÷8="
$t.Cell(17,5).Range.Text = "866÷5="
